# unify the conception of DataNode, DataTable, Entity.
# The "Property" worksheet is renamed to "DataNode" to match the unified
# naming scheme used across the other DataConfig workbooks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (was "Property") to the unified name "DataNode".
$ws.Name = "DataNode"

# Restore the cursor/selection position recorded in the saved file
# (frozen pane stays on row 9, but the active cell moved to D39).
$ws.Range("D39").Select()
